$wb = $excel.ActiveWorkbook

# This script applies updated market-board price/profit figures
# (scheduled data refresh) to each leve-profit worksheet.

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 156.35294
$ws.Range("I5").Value = 90.083336
$ws.Range("K5").Value = 90.083336
$ws.Range("M5").Value = 24.916664
$ws.Range("H9").Value = 55555656
$ws.Range("I9").Value = 200
$ws.Range("J9").Value = 111111110
$ws.Range("K9").Value = 200
$ws.Range("L9").Value = 111111110
$ws.Range("M9").Value = -31
$ws.Range("N9").Value = -111111448
$ws.Range("H18").Value = 2078.5789
$ws.Range("I18").Value = 2333.6
$ws.Range("K18").Value = 2333.6
$ws.Range("M18").Value = -2049.6
$ws.Range("H76").Value = 5498.5
$ws.Range("I76").Value = 4997
$ws.Range("K76").Value = 4997
$ws.Range("M76").Value = -4682
$ws.Range("H79").Value = 5498.5
$ws.Range("I79").Value = 4997
$ws.Range("K79").Value = 4997
$ws.Range("M79").Value = -3905
$ws.Range("H97").Value = 2033
$ws.Range("J97").Value = 1703.4445
$ws.Range("L97").Value = 5110.333500000001
$ws.Range("N97").Value = -6102.333500000001
$ws.Range("H112").Value = 20772.691
$ws.Range("I112").Value = 3749.6667
$ws.Range("K112").Value = 11249.0001
$ws.Range("M112").Value = -10141.0001
$ws.Range("H116").Value = 2017871.1
$ws.Range("J116").Value = 2499
$ws.Range("L116").Value = 2499
$ws.Range("N116").Value = -9383
$ws.Range("H132").Value = 5080.2
$ws.Range("I132").Value = 5130.4443
$ws.Range("J132").Value = 4628
$ws.Range("K132").Value = 15391.3329
$ws.Range("L132").Value = 13884
$ws.Range("M132").Value = -12861.3329
$ws.Range("N132").Value = -18944
$ws.Range("H138").Value = 4618.8853
$ws.Range("I138").Value = 987.0625
$ws.Range("J138").Value = 5910.2
$ws.Range("K138").Value = 2961.1875
$ws.Range("L138").Value = 17730.6
$ws.Range("M138").Value = 2178.8125
$ws.Range("N138").Value = -28010.6

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2919.625
$ws.Range("I32").Value = 2861.946
$ws.Range("K32").Value = 2861.946
$ws.Range("M32").Value = -2574.946
$ws.Range("H61").Value = 16969.555
$ws.Range("I61").Value = 41012.668
$ws.Range("J61").Value = 4948
$ws.Range("K61").Value = 41012.668
$ws.Range("L61").Value = 4948
$ws.Range("M61").Value = -40800.668
$ws.Range("N61").Value = -5372
$ws.Range("H74").Value = 132270.55
$ws.Range("I74").Value = 144997.6
$ws.Range("K74").Value = 144997.6
$ws.Range("M74").Value = -144123.6
$ws.Range("H77").Value = 132270.55
$ws.Range("I77").Value = 144997.6
$ws.Range("K77").Value = 724988
$ws.Range("M77").Value = -720620
$ws.Range("H136").Value = 16969.555
$ws.Range("I136").Value = 41012.668
$ws.Range("J136").Value = 4948
$ws.Range("K136").Value = 123038.004
$ws.Range("L136").Value = 14844
$ws.Range("M136").Value = -120488.004
$ws.Range("N136").Value = -19944

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 396
$ws.Range("I22").Value = 396
$ws.Range("K22").Value = 396
$ws.Range("M22").Value = -223
$ws.Range("H75").Value = 1000
$ws.Range("I75").Value = 1000
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 1000
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -64
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 1000
$ws.Range("I78").Value = 1000
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 3000
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = 1680
$ws.Range("N78").ClearContents()

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 5555693.5
$ws.Range("I7").Value = 153.90909
$ws.Range("J7").Value = 14285827
$ws.Range("K7").Value = 153.90909
$ws.Range("L7").Value = 14285827
$ws.Range("M7").Value = -40.90908999999999
$ws.Range("N7").Value = -14286053
$ws.Range("H31").Value = 1988.2319
$ws.Range("I31").Value = 1436.5807
$ws.Range("J31").Value = 2438.2632
$ws.Range("K31").Value = 1436.5807
$ws.Range("L31").Value = 2438.2632
$ws.Range("M31").Value = -1141.5807
$ws.Range("N31").Value = -3028.2632
$ws.Range("H34").Value = 1988.2319
$ws.Range("I34").Value = 1436.5807
$ws.Range("J34").Value = 2438.2632
$ws.Range("K34").Value = 1436.5807
$ws.Range("L34").Value = 2438.2632
$ws.Range("M34").Value = -1234.5807
$ws.Range("N34").Value = -2842.2632
$ws.Range("H58").Value = 6710.5483
$ws.Range("I58").Value = 9172.294
$ws.Range("J58").Value = 3721.2856
$ws.Range("K58").Value = 9172.294
$ws.Range("L58").Value = 3721.2856
$ws.Range("M58").Value = -8969.294
$ws.Range("N58").Value = -4127.2856
$ws.Range("H99").Value = 504342.3
$ws.Range("J99").Value = 8750
$ws.Range("L99").Value = 8750
$ws.Range("N99").Value = -11746
$ws.Range("H126").Value = 504342.3
$ws.Range("J126").Value = 8750
$ws.Range("L126").Value = 26250
$ws.Range("N126").Value = -31190
$ws.Range("H136").Value = 6710.5483
$ws.Range("I136").Value = 9172.294
$ws.Range("J136").Value = 3721.2856
$ws.Range("K136").Value = 27516.882
$ws.Range("L136").Value = 11163.8568
$ws.Range("M136").Value = -24966.882
$ws.Range("N136").Value = -16263.8568

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 313529.28
$ws.Range("J5").Value = 835116.5
$ws.Range("L5").Value = 2505349.5
$ws.Range("N5").Value = -2505573.5
$ws.Range("H8").Value = 359.3846
$ws.Range("I8").Value = 359.3846
$ws.Range("K8").Value = 1078.1538
$ws.Range("M8").Value = -939.1538
$ws.Range("H15").Value = 199
$ws.Range("J15").Value = 299
$ws.Range("L15").Value = 897
$ws.Range("N15").Value = -1177
$ws.Range("H32").Value = 1278.4286
$ws.Range("I32").Value = 450
$ws.Range("J32").Value = 1416.5
$ws.Range("K32").Value = 1350
$ws.Range("L32").Value = 4249.5
$ws.Range("M32").Value = -1067
$ws.Range("N32").Value = -4815.5
$ws.Range("H56").Value = 5151.1113
$ws.Range("I56").Value = 5151.1113
$ws.Range("K56").Value = 5151.1113
$ws.Range("M56").Value = -4621.1113
$ws.Range("H110").Value = 58163.5
$ws.Range("I110").Value = 49994
$ws.Range("K110").Value = 149982
$ws.Range("M110").Value = -145892
$ws.Range("H113").Value = 811.94446
$ws.Range("I113").Value = 552.6923
$ws.Range("K113").Value = 1658.0769
$ws.Range("M113").Value = 511.9231
$ws.Range("H121").Value = 373049.8
$ws.Range("J121").Value = 863.6429000000001
$ws.Range("L121").Value = 2590.9287
$ws.Range("N121").Value = -5210.9287
$ws.Range("H135").Value = 313529.28
$ws.Range("J135").Value = 835116.5
$ws.Range("L135").Value = 7516048.5
$ws.Range("N135").Value = -7521118.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 11890.591
$ws.Range("I122").Value = 12075.904
$ws.Range("K122").Value = 36227.712
$ws.Range("M122").Value = -33777.712

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5351.75
$ws.Range("I22").Value = 6759.875
$ws.Range("J22").Value = 3662
$ws.Range("K22").Value = 6759.875
$ws.Range("L22").Value = 3662
$ws.Range("M22").Value = -6464.875
$ws.Range("N22").Value = -4252
$ws.Range("H27").Value = 5351.75
$ws.Range("I27").Value = 6759.875
$ws.Range("J27").Value = 3662
$ws.Range("K27").Value = 6759.875
$ws.Range("L27").Value = 3662
$ws.Range("M27").Value = -6652.875
$ws.Range("N27").Value = -3876
$ws.Range("H40").Value = 56161
$ws.Range("I40").Value = 59417.31
$ws.Range("J40").Value = 34995
$ws.Range("K40").Value = 59417.31
$ws.Range("L40").Value = 34995
$ws.Range("M40").Value = -59281.31
$ws.Range("N40").Value = -35267
$ws.Range("H43").Value = 15127.75
$ws.Range("I43").Value = 15012
$ws.Range("J43").Value = 15166.333
$ws.Range("K43").Value = 15012
$ws.Range("L43").Value = 15166.333
$ws.Range("M43").Value = -14819
$ws.Range("N43").Value = -15552.333
$ws.Range("H61").Value = 1511.3636
$ws.Range("I61").Value = 1251.4286
$ws.Range("K61").Value = 1251.4286
$ws.Range("M61").Value = -1049.4286
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H82").Value = 2131.0688
$ws.Range("I82").Value = 2583.4285
$ws.Range("J82").Value = 1708.8667
$ws.Range("K82").Value = 2583.4285
$ws.Range("L82").Value = 1708.8667
$ws.Range("M82").Value = -2222.4285
$ws.Range("N82").Value = -2430.8667
$ws.Range("H85").Value = 2131.0688
$ws.Range("I85").Value = 2583.4285
$ws.Range("J85").Value = 1708.8667
$ws.Range("K85").Value = 2583.4285
$ws.Range("L85").Value = 1708.8667
$ws.Range("M85").Value = -1335.4285
$ws.Range("N85").Value = -4204.8667
$ws.Range("H113").Value = 1511.3636
$ws.Range("I113").Value = 1251.4286
$ws.Range("K113").Value = 1251.4286
$ws.Range("M113").Value = 918.5714

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 24397.8
$ws.Range("J37").Value = 26997.5
$ws.Range("L37").Value = 26997.5
$ws.Range("N37").Value = -27403.5
$ws.Range("H113").Value = 2988.9666
$ws.Range("I113").Value = 888.7895
$ws.Range("J113").Value = 6616.5454
$ws.Range("K113").Value = 2666.3685
$ws.Range("L113").Value = 19849.6362
$ws.Range("M113").Value = -496.3685
$ws.Range("N113").Value = -24189.6362
